$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1942446043165468
$ws.Range("C2").Value = 0.539568345323741
$ws.Range("J2").Value = 0.02158273381294964
$ws.Range("P2").Value = 0.1510791366906475
$ws.Range("S2").Value = 0.09352517985611511
$ws.Range("B3").Value = 0.03846153846153846
$ws.Range("C3").Value = 0.03846153846153846
$ws.Range("J3").Value = 0.02564102564102564
$ws.Range("P3").Value = 0.7564102564102564
$ws.Range("S3").Value = 0.141025641025641
$ws.Range("B6").Value = 0.07079646017699115
$ws.Range("F6").Value = 0.05309734513274336
$ws.Range("J6").Value = 0.2831858407079646
$ws.Range("O6").Value = 0.01769911504424779
$ws.Range("Q6").Value = 0.1946902654867257
$ws.Range("R6").Value = 0.07079646017699115
$ws.Range("S6").Value = 0.3097345132743363
$ws.Range("B7").Value = 0.1009174311926606
$ws.Range("D7").Value = 0.03669724770642202
$ws.Range("F7").Value = 0.03669724770642202
$ws.Range("J7").Value = 0.2018348623853211
$ws.Range("Q7").Value = 0.1559633027522936
$ws.Range("R7").Value = 0.04587155963302753
$ws.Range("S7").Value = 0.4220183486238532
$ws.Range("B8").Value = 0.06635071090047394
$ws.Range("D8").Value = 0.02369668246445497
$ws.Range("F8").Value = 0.06161137440758294
$ws.Range("J8").Value = 0.1042654028436019
$ws.Range("O8").Value = 0.02369668246445497
$ws.Range("Q8").Value = 0.2227488151658768
$ws.Range("R8").Value = 0.0995260663507109
$ws.Range("S8").Value = 0.3981042654028436
$ws.Range("B9").Value = 0.04301075268817205
$ws.Range("D9").Value = 0.01075268817204301
$ws.Range("F9").Value = 0.1075268817204301
$ws.Range("J9").Value = 0.1720430107526882
$ws.Range("O9").Value = 0.02150537634408602
$ws.Range("Q9").Value = 0.1397849462365591
$ws.Range("R9").Value = 0.1290322580645161
$ws.Range("S9").Value = 0.3763440860215054
$ws.Range("B10").Value = 0.115569823434992
$ws.Range("D10").Value = 0.01284109149277689
$ws.Range("E10").Value = 0.001605136436597111
$ws.Range("F10").Value = 0.06581059390048154
$ws.Range("J10").Value = 0.1107544141252006
$ws.Range("O10").Value = 0.01605136436597111
$ws.Range("Q10").Value = 0.1797752808988764
$ws.Range("R10").Value = 0.08667736757624397
$ws.Range("S10").Value = 0.4109149277688603
$ws.Range("G11").Value = 0.1314285714285714
$ws.Range("J11").Value = 0.07428571428571429
$ws.Range("K11").Value = 0.1771428571428571
$ws.Range("L11").Value = 0.5942857142857143
$ws.Range("S11").Value = 0.02285714285714286
$ws.Range("J12").Value = 0.2280701754385965
$ws.Range("K12").Value = 0.02631578947368421
$ws.Range("L12").Value = 0.04385964912280702
$ws.Range("S12").Value = 0.03508771929824561
$ws.Range("G13").Value = 0.5416666666666666
$ws.Range("J13").Value = 0.4166666666666667
$ws.Range("S13").Value = 0.04166666666666666
$ws.Range("F15").Value = 0.03773584905660377
$ws.Range("H15").Value = 0.1886792452830189
$ws.Range("I15").Value = 0.05660377358490566
$ws.Range("J15").Value = 0.3113207547169811
$ws.Range("K15").Value = 0.07547169811320754
$ws.Range("M15").Value = 0.009433962264150943
$ws.Range("O15").Value = 0.09433962264150944
$ws.Range("S15").Value = 0.2264150943396226
$ws.Range("F16").Value = 0.03488372093023256
$ws.Range("H16").Value = 0.1162790697674419
$ws.Range("I16").Value = 0.1627906976744186
$ws.Range("J16").Value = 0.3953488372093023
$ws.Range("K16").Value = 0.1162790697674419
$ws.Range("O16").Value = 0.05813953488372093
$ws.Range("S16").Value = 0.1162790697674419
$ws.Range("F17").Value = 0.02870813397129187
$ws.Range("H17").Value = 0.1818181818181818
$ws.Range("I17").Value = 0.1148325358851675
$ws.Range("J17").Value = 0.3779904306220095
$ws.Range("K17").Value = 0.1291866028708134
$ws.Range("M17").Value = 0.03827751196172249
$ws.Range("N17").Value = 0.004784688995215311
$ws.Range("O17").Value = 0.04784688995215311
$ws.Range("S17").Value = 0.07655502392344497
$ws.Range("F18").Value = 0.02
$ws.Range("H18").Value = 0.19
$ws.Range("I18").Value = 0.05
$ws.Range("J18").Value = 0.4
$ws.Range("K18").Value = 0.1
$ws.Range("M18").Value = 0.01
$ws.Range("O18").Value = 0.07000000000000001
$ws.Range("S18").Value = 0.16
$ws.Range("F19").Value = 0.01555209953343701
$ws.Range("H19").Value = 0.1990668740279938
$ws.Range("I19").Value = 0.06687402799377916
$ws.Range("J19").Value = 0.3608087091757387
$ws.Range("K19").Value = 0.135303265940902
$ws.Range("M19").Value = 0.02332814930015552
$ws.Range("N19").Value = 0.001555209953343701
$ws.Range("O19").Value = 0.06687402799377916
$ws.Range("S19").Value = 0.1306376360808709
